# Word COM-interop script implementing the commit's change:
#  - remove the big bold/underlined title paragraph ("NLP - Sentiment Analysis"),
#    which causes every subsequent paragraph to shift up by one position
#  - the document's trailing "_GoBack" bookmark (which Word places at the last
#    edited location) ends up right after the "These hands-on challenges..."
#    paragraph instead of in the now-trailing empty paragraph
#  - the numbered list's level-1 format switches from "1." (decimal) to
#    "a)" (lowerLetter)

$d = $word.ActiveDocument

# 1) Delete the first paragraph (the document title). Using Range.Delete()
#    removes the paragraph mark along with its text, so every following
#    paragraph shifts up by one - exactly matching the diff, which drops the
#    title paragraph and reindexes the rest.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Delete()

# 2) Relocate the "_GoBack" bookmark so it sits immediately after the text of
#    the paragraph that now reads "These hands-on challenges sharpened my
#    production mindset and revealed key insights:" (a zero-length bookmark,
#    bookmarkStart immediately followed by bookmarkEnd, right after the run).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*These hands-on challenges*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Exclude the trailing paragraph mark from the range, landing on the
    # last real character of the paragraph's text.
    [void]$r.MoveEnd(1, -1)
    # Adding with a (non-empty) one-character-shorter range first, then
    # collapsing it in place to zero length, reliably places the bookmark
    # right after the run; adding it directly with a zero-length range
    # mis-positions it, so that two-step dance is intentional. This also
    # replaces/removes the old "_GoBack" bookmark, since bookmark names
    # must be unique within the document.
    $bm = $d.Bookmarks.Add("_GoBack", $r)
    $bmRange = $bm.Range
    $bmRange.Start = $bmRange.End
}

# 3) Change the single-level numbered list's formatting from "1." to "a)".
foreach ($p in $d.Paragraphs) {
    if ($p.Range.ListFormat.ListType -ne 0) {
        $lvl = $p.Range.ListFormat.ListTemplate.ListLevels.Item(1)
        $lvl.NumberFormat = "%1)"
        $lvl.NumberStyle = 4  # wdListNumberStyleLowercaseLetter
        break
    }
}
